$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix photosynthesis-related default input values
$ws.Range("B5").Value = 0.055
$ws.Range("B8").Value = 0.35

# Update the saved view: scroll back to top and select B5 instead of A35
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B5").Select()
